# Insert a new row at position 607 (shifts existing rows 607:652 down to 608:653)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(607).Insert()

# Populate the newly inserted row 607 with the new record's data
$ws.Cells.Item(607, 1).Value = 5
$ws.Cells.Item(607, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(607, 3).Value = "Maule"
$ws.Cells.Item(607, 4).Value = 44783
$ws.Cells.Item(607, 5).Value = 7
$ws.Cells.Item(607, 6).Value = 100112004
$ws.Cells.Item(607, 7).Value = "Cebolla"
$ws.Cells.Item(607, 8).Value = "Sin especificar"
$ws.Cells.Item(607, 9).Value = "1a (guarda)"
$ws.Cells.Item(607, 10).Value = 2500
$ws.Cells.Item(607, 11).Value = 10000
$ws.Cells.Item(607, 12).Value = 10000
$ws.Cells.Item(607, 13).Value = 10000
$ws.Cells.Item(607, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(607, 15).Value = "Región del Maule"
$ws.Cells.Item(607, 16).Value = 400
$ws.Cells.Item(607, 17).Value = 25
$ws.Cells.Item(607, 18).Value = "Hortaliza"
